# Apply fix: correct property_category values for 建物 (Building) and 汽車 (Car)
# sheets, which were incorrectly carrying the "land" value copied from the
# 土地 (Land) sheet, plus reconcile the row-index bookkeeping columns that
# shift because two new category rows ("building"/"car") were inserted
# upstream (commit "#5: property aircraft done").

$wb = $excel.ActiveWorkbook

# --- Sheet 2: 建物 (Building) -------------------------------------------
# Column I = property_category. Fix "land" -> "building" for data rows 2-6.
$wsBuilding = $wb.Worksheets.Item(2)
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"
$wsBuilding.Range("I4").Value = "building"
$wsBuilding.Range("I5").Value = "building"
$wsBuilding.Range("I6").Value = "building"

# --- Sheet 3: 汽車 (Car) --------------------------------------------------
# Column H = property_category. Fix "land" -> "car" for data rows 2-3.
$wsCar = $wb.Worksheets.Item(3)
$wsCar.Range("H2").Value = "car"
$wsCar.Range("H3").Value = "car"

# Row bookkeeping columns (A = index, N = total-row reference) shift down by
# 2 because of the two newly-inserted category rows upstream.
$wsCar.Range("A2").Value = 38
$wsCar.Range("N2").Value = 38
$wsCar.Range("A3").Value = 39
$wsCar.Range("N3").Value = 39

# --- Sheet 4: 存款 (Deposit) ----------------------------------------------
$wsDeposit = $wb.Worksheets.Item(4)
$wsDeposit.Range("A2").Value = 53

# --- Sheet 5: 具有相當價值之財產 (Valuable property) -----------------------
$wsValuable = $wb.Worksheets.Item(5)
$wsValuable.Range("A2").Value = 80
$wsValuable.Range("A3").Value = 81
$wsValuable.Range("A4").Value = 82

# --- Sheet 6: 債務 (Debt) --------------------------------------------------
$wsDebt = $wb.Worksheets.Item(6)
$wsDebt.Range("A2").Value = 92
$wsDebt.Range("A3").Value = 93
$wsDebt.Range("A4").Value = 94
